$p = $ppt.ActivePresentation

# --- Slide 1: reposition/resize "TextBox 103" (the "C" label) ---
# Target EMU values (from OOXML): off x=5186528 y=3212068, ext cx=354806 cy=369332
# Shape.Left/Top/Width/Height round-trip through single-precision point
# values, so the literals below are chosen to land exactly on the target
# EMUs once converted back on save (verified empirically against the host).
$s1 = $p.Slides.Item(1)
$shpC = $s1.Shapes.Item(45)
$shpC.Left = 408.3881
$shpC.Top = 252.9188
$shpC.Width = 27.9375
$shpC.Height = 29.0813

# --- Slide 4: merge the two runs "1- " and "P_2" into a single run ---
# The concatenated text is already "1- P_2", so assigning that same string
# directly keeps (merges into) the existing two-run structure. Round-trip
# through an unrelated placeholder value first to force a full rewrite,
# collapsing the paragraph down to a single run before setting the final text.
$s4 = $p.Slides.Item(4)
$shpLabel = $s4.Shapes.Item(48)
$shpLabel.TextFrame.TextRange.Text = "TEMP"
$shpLabel.TextFrame.TextRange.Text = "1- P_2"
